$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "Flow 6 Addresst"
$ws.Range("B18").Value = "Flow 6 City"
$ws.Range("C18").Value = " "
$ws.Range("D18").Value = "Flow 6 First Name"
$ws.Range("E18").Value = "Flow 7 Lasst Name"
$ws.Range("F18").Value = " "
$ws.Range("G18").Value = "'666"
$ws.Range("H18").Value = "'6666"

$ws.Range("A19").Value = "Address Flow 1"
$ws.Range("B19").Value = "City Flow 1"
$ws.Range("C19").Value = " "
$ws.Range("D19").Value = "first flow 1"
$ws.Range("E19").Value = "last flow 2"
$ws.Range("F19").Value = " "
$ws.Range("G19").Value = "'111"
$ws.Range("H19").Value = "'1111"

$ws.Range("A20").Value = "Adrress 25 f1 @#`$%^&*!(#)#*"
$ws.Range("B20").Value = "25 F1 City"
$ws.Range("C20").Value = " "
$ws.Range("D20").Value = "25 f1 first name"
$ws.Range("E20").Value = "25 f1 last name"
$ws.Range("F20").Value = "'25012334567955"
$ws.Range("G20").Value = "'251"
$ws.Range("H20").Value = "'2501"

$ws.Range("A21").Value = "address f5"
$ws.Range("B21").Value = "city f5"
$ws.Range("C21").Value = " "
$ws.Range("D21").Value = "first f5"
$ws.Range("E21").Value = "last f5"
$ws.Range("F21").Value = " "
$ws.Range("G21").Value = "'555"
$ws.Range("H21").Value = "'5555"

$ws.Range("A22").Value = "goo"
$ws.Range("B22").Value = "ho chi minh"
$ws.Range("C22").Value = " "
$ws.Range("D22").Value = "firs"
$ws.Range("E22").Value = "last"
$ws.Range("F22").Value = " "
$ws.Range("G22").Value = "'5566"
$ws.Range("H22").Value = "'66"

$ws.Range("A23").Value = "sdfsdf"
$ws.Range("B23").Value = "gsdfgd"
$ws.Range("C23").Value = " "
$ws.Range("D23").Value = "df"
$ws.Range("E23").Value = "fsd"
$ws.Range("F23").Value = " "
$ws.Range("G23").Value = "'444"
$ws.Range("H23").Value = "'333"

$ws.Range("A24").Value = "Address Flow 4"
$ws.Range("B24").Value = "City Flow 4"
$ws.Range("C24").Value = " "
$ws.Range("D24").Value = "Flow 4 First"
$ws.Range("E24").Value = "Flow 4 Last"
$ws.Range("F24").Value = " "
$ws.Range("G24").Value = "'444"
$ws.Range("H24").Value = "'4444"

$ws.Range("A25").Value = "Adress flow 2"
$ws.Range("B25").Value = "City Flow 2"
$ws.Range("C25").Value = " "
$ws.Range("D25").Value = "first"
$ws.Range("E25").Value = "last flow 2"
$ws.Range("F25").Value = " "
$ws.Range("G25").Value = "'222"
$ws.Range("H25").Value = "'2222"

$ws.Range("A26").Value = "sdafsd"
$ws.Range("B26").Value = "dgdfg"
$ws.Range("C26").Value = " "
$ws.Range("D26").Value = "first"
$ws.Range("E26").Value = "ggg"
$ws.Range("F26").Value = "'095756756757656"
$ws.Range("G26").Value = "'444"
$ws.Range("H26").Value = "fsdfsd"

$ws.Range("A27").Value = "25 f4 Address"
$ws.Range("B27").Value = "25 f4 city"
$ws.Range("C27").Value = " "
$ws.Range("D27").Value = "25 f4 first"
$ws.Range("E27").Value = "25 f4 last"
$ws.Range("F27").Value = " "
$ws.Range("G27").Value = "'2544"
$ws.Range("H27").Value = "{{address}}"

$ws.Range("A28").Value = "hhh"
$ws.Range("B28").Value = "hhh"
$ws.Range("C28").Value = " "
$ws.Range("D28").Value = "gg"
$ws.Range("E28").Value = "hhh"
$ws.Range("F28").Value = " "
$ws.Range("G28").Value = "'777"
$ws.Range("H28").Value = "'777"

$ws.Range("A29").Value = "f3 address"
$ws.Range("B29").Value = "f2 city"
$ws.Range("C29").Value = " "
$ws.Range("D29").Value = "f2. first"
$ws.Range("E29").Value = "f2 last"
$ws.Range("F29").Value = " "
$ws.Range("G29").Value = "'222"
$ws.Range("H29").Value = "'2222"

$ws.Range("A30").Value = "Adress Flow 3 ( Step 2)"
$ws.Range("B30").Value = "City Flow 3 ( Step 5@)"
$ws.Range("C30").Value = " "
$ws.Range("D30").Value = "First Flow 3"
$ws.Range("E30").Value = "last Flow 3 ( Step 5)"
$ws.Range("F30").Value = " "
$ws.Range("G30").Value = "'333"
$ws.Range("H30").Value = "'3333"

$ws.Range("A31").Value = "25 f2 Address (2)"
$ws.Range("B31").Value = "25 f2 City ( 3)"
$ws.Range("C31").Value = " "
$ws.Range("D31").Value = "25 f2 First (1)"
$ws.Range("E31").Value = "25 F2 Last (4)"
$ws.Range("F31").Value = "'2502834949444"
$ws.Range("G31").Value = "'252"
$ws.Range("H31").Value = "'250205"

$ws.Range("A32").Value = "f3 second address"
$ws.Range("B32").Value = "f3 4 city"
$ws.Range("C32").Value = " "
$ws.Range("D32").Value = "f3 first (s1)"
$ws.Range("E32").Value = "f3 five last"
$ws.Range("F32").Value = " "
$ws.Range("G32").Value = "'333"
$ws.Range("H32").Value = "'3333"

$ws.Range("A33").Value = "f5 address 25"
$ws.Range("B33").Value = "f5 city 25"
$ws.Range("C33").Value = " "
$ws.Range("D33").Value = "f5 first 25"
$ws.Range("E33").Value = "f5 last 25"
$ws.Range("F33").Value = "{{credit_debit_number}}"
$ws.Range("G33").Value = "'2505"
$ws.Range("H33").Value = "'250505"

